$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 80; this pushes existing rows 80..188 down to 81..189
$ws.Rows.Item(80).Insert()

# Populate the newly inserted row 80 with the new data entry.
$ws.Cells.Item(80, 1).Value = 4
$ws.Cells.Item(80, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(80, 3).Value = "Los Lagos"
$ws.Cells.Item(80, 4).Value = 44579
$ws.Cells.Item(80, 4).NumberFormat = $ws.Cells.Item(81, 4).NumberFormat
$ws.Cells.Item(80, 5).Value = 10
$ws.Cells.Item(80, 6).Value = 100112017
$ws.Cells.Item(80, 7).Value = "Apio"
$ws.Cells.Item(80, 8).Value = "Americana (o)"
$ws.Cells.Item(80, 9).Value = "Primera"
$ws.Cells.Item(80, 10).Value = 50
$ws.Cells.Item(80, 11).Value = 11000
$ws.Cells.Item(80, 12).Value = 12000
$ws.Cells.Item(80, 13).Value = 11500
$ws.Cells.Item(80, 14).Value = "$/docena de matas"
$ws.Cells.Item(80, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(80, 16).Value = 1917
$ws.Cells.Item(80, 17).Value = 6
$ws.Cells.Item(80, 18).Value = "Hortaliza"
